$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.57
$ws.Range("C2").Value = 1.57
$ws.Range("D2").Value = 1.57
$ws.Range("E2").Value = 1.51
$ws.Range("F2").Value = 1.63
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 1.51
$ws.Range("I2").Value = 1.57
$ws.Range("J2").Value = 1.54
$ws.Range("K2").Value = 1.6
$ws.Range("L2").Value = 1.57
$ws.Range("M2").Value = 1.57
$ws.Range("N2").Value = 1.51
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 1.57
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.63
$ws.Range("T2").Value = 1.6
$ws.Range("U2").Value = 1.63
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.7
$ws.Range("X2").Value = 1.63
$ws.Range("Y2").Value = 1.57
$ws.Range("Z2").Value = 1.63
$ws.Range("AA2").Value = 1.57
$ws.Range("AB2").Value = 1.57
$ws.Range("AC2").Value = 1.7
$ws.Range("AD2").Value = 1.63
$ws.Range("AE2").Value = 1.57
$ws.Range("AF2").Value = 1.57
$ws.Range("AG2").Value = 1.57
$ws.Range("B3").Value = 1.4
$ws.Range("C3").Value = 1.4
$ws.Range("D3").Value = 1.39
$ws.Range("E3").Value = 1.36
$ws.Range("F3").Value = 1.43
$ws.Range("G3").Value = 1.4
$ws.Range("H3").Value = 1.36
$ws.Range("I3").Value = 1.41
$ws.Range("J3").Value = 1.38
$ws.Range("K3").Value = 1.43
$ws.Range("L3").Value = 1.4
$ws.Range("M3").Value = 1.39
$ws.Range("N3").Value = 1.34
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 1.4
$ws.Range("Q3").Value = 1.46
$ws.Range("R3").Value = 1.38
$ws.Range("S3").Value = 1.43
$ws.Range("T3").Value = 1.41
$ws.Range("U3").Value = 1.43
$ws.Range("V3").Value = 1.39
$ws.Range("W3").Value = 1.46
$ws.Range("X3").Value = 1.43
$ws.Range("Y3").Value = 1.4
$ws.Range("Z3").Value = 1.43
$ws.Range("AA3").Value = 1.39
$ws.Range("AB3").Value = 1.4
$ws.Range("AC3").Value = 1.48
$ws.Range("AD3").Value = 1.44
$ws.Range("AE3").Value = 1.39
$ws.Range("AF3").Value = 1.4
$ws.Range("AG3").Value = 1.4
$ws.Range("B4").Value = 1.36
$ws.Range("C4").Value = 1.35
$ws.Range("D4").Value = 1.35
$ws.Range("E4").Value = 1.33
$ws.Range("F4").Value = 1.38
$ws.Range("G4").Value = 1.36
$ws.Range("H4").Value = 1.33
$ws.Range("I4").Value = 1.37
$ws.Range("J4").Value = 1.34
$ws.Range("K4").Value = 1.38
$ws.Range("L4").Value = 1.35
$ws.Range("M4").Value = 1.35
$ws.Range("N4").Value = 1.31
$ws.Range("O4").Value = 1.37
$ws.Range("P4").Value = 1.37
$ws.Range("Q4").Value = 1.41
$ws.Range("R4").Value = 1.34
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 1.38
$ws.Range("U4").Value = 1.39
$ws.Range("V4").Value = 1.35
$ws.Range("W4").Value = 1.42
$ws.Range("X4").Value = 1.39
$ws.Range("Y4").Value = 1.36
$ws.Range("Z4").Value = 1.38
$ws.Range("AA4").Value = 1.35
$ws.Range("AB4").Value = 1.36
$ws.Range("AC4").Value = 1.43
$ws.Range("AD4").Value = 1.39
$ws.Range("AE4").Value = 1.35
$ws.Range("AF4").Value = 1.35
$ws.Range("AG4").Value = 1.36
$ws.Range("B5").Value = 1.34
$ws.Range("C5").Value = 1.33
$ws.Range("D5").Value = 1.33
$ws.Range("E5").Value = 1.31
$ws.Range("F5").Value = 1.36
$ws.Range("G5").Value = 1.34
$ws.Range("H5").Value = 1.31
$ws.Range("I5").Value = 1.35
$ws.Range("J5").Value = 1.32
$ws.Range("K5").Value = 1.36
$ws.Range("L5").Value = 1.33
$ws.Range("M5").Value = 1.34
$ws.Range("N5").Value = 1.29
$ws.Range("O5").Value = 1.35
$ws.Range("P5").Value = 1.35
$ws.Range("Q5").Value = 1.39
$ws.Range("R5").Value = 1.32
$ws.Range("S5").Value = 1.37
$ws.Range("T5").Value = 1.36
$ws.Range("U5").Value = 1.37
$ws.Range("V5").Value = 1.33
$ws.Range("W5").Value = 1.4
$ws.Range("X5").Value = 1.37
$ws.Range("Y5").Value = 1.34
$ws.Range("Z5").Value = 1.36
$ws.Range("AA5").Value = 1.33
$ws.Range("AB5").Value = 1.34
$ws.Range("AC5").Value = 1.41
$ws.Range("AD5").Value = 1.37
$ws.Range("AE5").Value = 1.33
$ws.Range("AF5").Value = 1.34
$ws.Range("AG5").Value = 1.34
$ws.Range("B6").Value = 1.33
$ws.Range("C6").Value = 1.32
$ws.Range("D6").Value = 1.32
$ws.Range("E6").Value = 1.3
$ws.Range("F6").Value = 1.35
$ws.Range("G6").Value = 1.32
$ws.Range("H6").Value = 1.3
$ws.Range("I6").Value = 1.34
$ws.Range("J6").Value = 1.32
$ws.Range("K6").Value = 1.35
$ws.Range("L6").Value = 1.32
$ws.Range("M6").Value = 1.33
$ws.Range("N6").Value = 1.28
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 1.34
$ws.Range("Q6").Value = 1.38
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 1.35
$ws.Range("U6").Value = 1.36
$ws.Range("V6").Value = 1.32
$ws.Range("W6").Value = 1.39
$ws.Range("X6").Value = 1.36
$ws.Range("Y6").Value = 1.33
$ws.Range("Z6").Value = 1.35
$ws.Range("AA6").Value = 1.32
$ws.Range("AB6").Value = 1.33
$ws.Range("AC6").Value = 1.4
$ws.Range("AD6").Value = 1.36
$ws.Range("AE6").Value = 1.33
$ws.Range("AF6").Value = 1.33
$ws.Range("AG6").Value = 1.33
$ws.Range("B7").Value = 1.33
$ws.Range("C7").Value = 1.32
$ws.Range("D7").Value = 1.32
$ws.Range("E7").Value = 1.3
$ws.Range("F7").Value = 1.35
$ws.Range("G7").Value = 1.32
$ws.Range("H7").Value = 1.3
$ws.Range("I7").Value = 1.33
$ws.Range("J7").Value = 1.31
$ws.Range("K7").Value = 1.35
$ws.Range("L7").Value = 1.32
$ws.Range("M7").Value = 1.32
$ws.Range("N7").Value = 1.28
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 1.33
$ws.Range("Q7").Value = 1.37
$ws.Range("R7").Value = 1.31
$ws.Range("S7").Value = 1.36
$ws.Range("T7").Value = 1.35
$ws.Range("U7").Value = 1.35
$ws.Range("V7").Value = 1.32
$ws.Range("W7").Value = 1.39
$ws.Range("X7").Value = 1.36
$ws.Range("Y7").Value = 1.33
$ws.Range("Z7").Value = 1.34
$ws.Range("AA7").Value = 1.32
$ws.Range("AB7").Value = 1.33
$ws.Range("AC7").Value = 1.39
$ws.Range("AD7").Value = 1.36
$ws.Range("AE7").Value = 1.32
$ws.Range("AF7").Value = 1.32
$ws.Range("AG7").Value = 1.33
$ws.Range("B8").Value = 1.32
$ws.Range("C8").Value = 1.32
$ws.Range("D8").Value = 1.32
$ws.Range("E8").Value = 1.3
$ws.Range("F8").Value = 1.35
$ws.Range("G8").Value = 1.32
$ws.Range("H8").Value = 1.3
$ws.Range("I8").Value = 1.33
$ws.Range("J8").Value = 1.31
$ws.Range("K8").Value = 1.35
$ws.Range("L8").Value = 1.32
$ws.Range("M8").Value = 1.32
$ws.Range("N8").Value = 1.28
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.33
$ws.Range("Q8").Value = 1.36
$ws.Range("R8").Value = 1.31
$ws.Range("S8").Value = 1.35
$ws.Range("T8").Value = 1.35
$ws.Range("U8").Value = 1.35
$ws.Range("V8").Value = 1.32
$ws.Range("W8").Value = 1.61
$ws.Range("X8").Value = 1.36
$ws.Range("Y8").Value = 1.33
$ws.Range("Z8").Value = 1.34
$ws.Range("AA8").Value = 1.32
$ws.Range("AB8").Value = 1.33
$ws.Range("AC8").Value = 1.39
$ws.Range("AD8").Value = 1.36
$ws.Range("AE8").Value = 1.32
$ws.Range("AF8").Value = 1.32
$ws.Range("AG8").Value = 1.33
$ws.Range("B9").Value = 1.33
$ws.Range("C9").Value = 1.33
$ws.Range("D9").Value = 1.33
$ws.Range("E9").Value = 1.31
$ws.Range("F9").Value = 1.35
$ws.Range("G9").Value = 1.33
$ws.Range("H9").Value = 1.31
$ws.Range("I9").Value = 1.34
$ws.Range("J9").Value = 1.32
$ws.Range("K9").Value = 1.35
$ws.Range("L9").Value = 1.32
$ws.Range("M9").Value = 1.33
$ws.Range("N9").Value = 1.29
$ws.Range("O9").Value = 1.34
$ws.Range("P9").Value = 1.34
$ws.Range("Q9").Value = 1.37
$ws.Range("R9").Value = 1.32
$ws.Range("S9").Value = 1.36
$ws.Range("T9").Value = 1.35
$ws.Range("U9").Value = 1.36
$ws.Range("V9").Value = 1.32
$ws.Range("W9").Value = 2
$ws.Range("X9").Value = 1.37
$ws.Range("Y9").Value = 1.33
$ws.Range("Z9").Value = 1.35
$ws.Range("AA9").Value = 1.32
$ws.Range("AB9").Value = 1.33
$ws.Range("AC9").Value = 1.4
$ws.Range("AD9").Value = 1.38
$ws.Range("AE9").Value = 1.33
$ws.Range("AF9").Value = 1.33
$ws.Range("AG9").Value = 1.33
